$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = @"
Pipeline(steps=[('scaler', RobustScaler()),
                ('selector',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1',
                                                     random_state=42))),
                ('model',
                 AdaBoostClassifier(estimator=RandomForestClassifier(class_weight='balanced',
                                                                     max_depth=5,
                                                                     max_features='log2',
                                                                     min_samples_leaf=5,
                                                                     n_estimators=50,
                                                                     random_state=42),
                                    n_estimators=10, random_state=42))])
"@
$ws.Range("B2").Value = 0.657142857142857
$ws.Range("C2").Value = @"
{'selector': SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1', random_state=42)), 'scaler': RobustScaler(), 'model__n_estimators': 10, 'model__estimator__n_estimators': 50, 'model__estimator__min_samples_split': 2, 'model__estimator__min_samples_leaf': 5, 'model__estimator__max_features': 'log2', 'model__estimator__max_depth': 5, 'model__estimator__class_weight': 'balanced'}
"@
$ws.Range("D2").Value = 0.4615384615384615
$ws.Range("E2").Value = @"
[1 1 0 0 1 0 0 0 0 1 0 1]
"@
$ws.Range("F2").Value = @"
[0 0 0 1 1 1 0 1 1 1 1 1]
"@
$ws.Range("G2").Value = 77
$ws.Range("H2").Value = 0.9696503918022907
$ws.Range("I2").Value = 0.008136686480157019
$ws.Range("J2").Value = 0.5723930078360459
$ws.Range("K2").Value = 0.05434431720975728
$ws.Rows.Item(2).AutoFit()

$ws.Range("A3").Value = @"
Pipeline(steps=[('scaler', StandardScaler()), ('selector', None),
                ('model',
                 AdaBoostClassifier(estimator=RandomForestClassifier(max_depth=1,
                                                                     max_features='log2',
                                                                     min_samples_leaf=5,
                                                                     min_samples_split=4,
                                                                     n_estimators=5,
                                                                     random_state=42),
                                    n_estimators=10, random_state=42))])
"@
$ws.Range("B3").Value = 0.6
$ws.Range("C3").Value = @"
{'selector': None, 'scaler': StandardScaler(), 'model__n_estimators': 10, 'model__estimator__n_estimators': 5, 'model__estimator__min_samples_split': 4, 'model__estimator__min_samples_leaf': 5, 'model__estimator__max_features': 'log2', 'model__estimator__max_depth': 1, 'model__estimator__class_weight': None}
"@
$ws.Range("D3").Value = 0.6666666666666665
$ws.Range("E3").Value = @"
[1 1 0 1 0 0 1 0 1 1 1 0]
"@
$ws.Range("F3").Value = @"
[0 1 1 1 1 1 1 1 1 1 1 1]
"@
$ws.Range("G3").Value = 69
$ws.Range("H3").Value = 0.9782581453634085
$ws.Range("I3").Value = 0.007066276457252609
$ws.Range("J3").Value = 0.5020050125313283
$ws.Range("K3").Value = 0.0885618519312343
$ws.Rows.Item(3).AutoFit()

$ws.Range("A4").Value = @"
Pipeline(steps=[('scaler', RobustScaler()), ('selector', None),
                ('model',
                 AdaBoostClassifier(estimator=RandomForestClassifier(max_depth=2,
                                                                     min_samples_split=3,
                                                                     n_estimators=5,
                                                                     random_state=42),
                                    n_estimators=5, random_state=42))])
"@
$ws.Range("B4").Value = 0.6571428571428573
$ws.Range("C4").Value = @"
{'selector': None, 'scaler': RobustScaler(), 'model__n_estimators': 5, 'model__estimator__n_estimators': 5, 'model__estimator__min_samples_split': 3, 'model__estimator__min_samples_leaf': 1, 'model__estimator__max_features': 'sqrt', 'model__estimator__max_depth': 2, 'model__estimator__class_weight': None}
"@
$ws.Range("D4").Value = 0.823529411764706
$ws.Range("E4").Value = @"
[1 0 1 1 1 1 0 1 0 1 0 1]
"@
$ws.Range("F4").Value = @"
[1 0 1 1 1 1 0 0 1 1 1 1]
"@
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 0.9731481481481481
$ws.Range("I4").Value = 0.006601313070372444
$ws.Range("J4").Value = 0.5294973544973545
$ws.Range("K4").Value = 0.08543319076137006
$ws.Rows.Item(4).AutoFit()

$ws.Range("A5").Value = @"
Pipeline(steps=[('scaler', None),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f3a64631970>),
                ('model',
                 AdaBoostClassifier(estimator=RandomForestClassifier(max_depth=5,
                                                                     min_samples_leaf=3,
                                                                     min_samples_split=5,
                                                                     n_estimators=50,
                                                                     random_state=42),
                                    n_estimators=5, random_state=42))])
"@
$ws.Range("B5").Value = 0.6095238095238095
$ws.Range("C5").Value = @"
{'selector': <__main__.NamedFeatureSelector object at 0x7f3a6c09d0a0>, 'scaler': None, 'model__n_estimators': 5, 'model__estimator__n_estimators': 50, 'model__estimator__min_samples_split': 5, 'model__estimator__min_samples_leaf': 3, 'model__estimator__max_features': 'sqrt', 'model__estimator__max_depth': 5, 'model__estimator__class_weight': None}
"@
$ws.Range("D5").Value = 0.6666666666666666
$ws.Range("E5").Value = @"
[1 1 0 0 0 0 1 0 1 1 1 1]
"@
$ws.Range("F5").Value = @"
[1 0 0 1 0 1 1 1 1 0 1 1]
"@
$ws.Range("G5").Value = 11
$ws.Range("H5").Value = 0.9845972957084069
$ws.Range("I5").Value = 0.005650381380085078
$ws.Range("J5").Value = 0.5027630805408583
$ws.Range("K5").Value = 0.06999311855864526
$ws.Rows.Item(5).AutoFit()

$ws.Range("A6").Value = @"
Pipeline(steps=[('scaler', None),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f3a647d0670>),
                ('model',
                 AdaBoostClassifier(estimator=RandomForestClassifier(max_depth=3,
                                                                     max_features='log2',
                                                                     min_samples_leaf=4,
                                                                     min_samples_split=5,
                                                                     n_estimators=50,
                                                                     random_state=42),
                                    n_estimators=5, random_state=42))])
"@
$ws.Range("B6").Value = 0.6
$ws.Range("C6").Value = @"
{'selector': <__main__.NamedFeatureSelector object at 0x7f3a6c0afb50>, 'scaler': None, 'model__n_estimators': 5, 'model__estimator__n_estimators': 50, 'model__estimator__min_samples_split': 5, 'model__estimator__min_samples_leaf': 4, 'model__estimator__max_features': 'log2', 'model__estimator__max_depth': 3, 'model__estimator__class_weight': None}
"@
$ws.Range("D6").Value = 0.7142857142857143
$ws.Range("E6").Value = @"
[1 1 1 1 0 0 0 0 1 1 0 0]
"@
$ws.Range("F6").Value = @"
[1 1 1 0 0 1 0 1 1 1 1 0]
"@
$ws.Range("G6").Value = 14
$ws.Range("H6").Value = 0.9781599059376838
$ws.Range("I6").Value = 0.006286474816145809
$ws.Range("J6").Value = 0.5279247501469724
$ws.Range("K6").Value = 0.0874152984061609
$ws.Rows.Item(6).AutoFit()

Write-Host "done"
